# Generate Report for Handoff
#
# Re-running the localization status report bumps the "Latest Handoff
# Datetime" timestamps for the 2f0b77d5-e1a6-4cef-8301-b0682e9376c7.md
# handoff (rows 8, 9, 10, 12, 13, 14 on each language sheet; row 11 is a
# different file and is left untouched) and marks those same rows' new
# "Priority" column as "ht" on the per-language detail sheets.

$wb = $excel.ActiveWorkbook

$rows = @(8, 9, 10, 12, 13, 14)

# --- Overview sheet: column G = "Latest HO Xliff Generate Date" ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-08-29 12:22:24"
}

# --- zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime" ---
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-29 12:22:18"
}

# --- de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime" ---
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-29 12:22:24"
}
